$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 2.76
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 1.45
$ws.Range("O2").Value = 1.15
$ws.Range("Q2").Value = 2.7
$ws.Range("R2").Value = 1.15
$ws.Range("S2").Value = 6.2
$ws.Range("T2").Value = 2.22
$ws.Range("U2").Value = 1.53
$ws.Range("V2").Value = 1.17
$ws.Range("W2").Value = 1.94
$ws.Range("X2").Value = 8.8
$ws.Range("Y2").Value = 14.5
$ws.Range("Z2").Value = 48
$ws.Range("AA2").Value = 220
$ws.Range("AB2").Value = 6.6
$ws.Range("AC2").Value = 8.8
$ws.Range("AD2").Value = 29
$ws.Range("AE2").Value = 150
$ws.Range("AF2").Value = 12
$ws.Range("AG2").Value = 13.5
$ws.Range("AH2").Value = 36
$ws.Range("AI2").Value = 180
$ws.Range("AJ2").Value = 29
$ws.Range("AK2").Value = 36
$ws.Range("AL2").Value = 85
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 34
$ws.Range("AO2").Value = 280

# Row 3
$ws.Range("F3").Value = 2.54
$ws.Range("G3").Value = 2.68
$ws.Range("H3").Value = 3.3
$ws.Range("J3").Value = 2.86
$ws.Range("K3").Value = 3.05
$ws.Range("L3").Value = 1.49
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 2.54
$ws.Range("O3").Value = 1.57
$ws.Range("Q3").Value = 2.38
$ws.Range("R3").Value = 1.18
$ws.Range("S3").Value = 4.8
$ws.Range("T3").Value = 1.94
$ws.Range("U3").Value = 1.62
$ws.Range("V3").Value = 1.38
$ws.Range("W3").Value = 1.59
$ws.Range("X3").Value = 10.5
$ws.Range("Y3").Value = 11.5
$ws.Range("Z3").Value = 27
$ws.Range("AA3").Value = 90
$ws.Range("AB3").Value = 9
$ws.Range("AC3").Value = 8.8
$ws.Range("AD3").Value = 19
$ws.Range("AE3").Value = 70
$ws.Range("AF3").Value = 18.5
$ws.Range("AG3").Value = 15.5
$ws.Range("AH3").Value = 30
$ws.Range("AI3").Value = 100
$ws.Range("AJ3").Value = 55
$ws.Range("AK3").Value = 46
$ws.Range("AL3").Value = 70
$ws.Range("AM3").Value = 250
$ws.Range("AN3").Value = 55
$ws.Range("AO3").Value = 95

# Row 4
$ws.Range("F4").Value = 1.25
$ws.Range("G4").Value = 1.93
$ws.Range("H4").Value = 4.7
$ws.Range("I4").Value = 7.2
$ws.Range("J4").Value = 3.7
$ws.Range("K4").Value = 980
$ws.Range("L4").Value = 1.01
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 1.8
$ws.Range("O4").Value = 1.31
$ws.Range("P4").Value = 1.8
$ws.Range("Q4").Value = 1.81
$ws.Range("R4").Value = 1.27
$ws.Range("S4").Value = 2.98
$ws.Range("T4").Value = 1.01
$ws.Range("U4").Value = 1.01
$ws.Range("V4").Value = 1.01
$ws.Range("W4").Value = 2.06
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

# Row 5
$ws.Range("M5").Value = 1.02
$ws.Range("O5").Value = 1.12
$ws.Range("Q5").Value = 1.36

# Row 7
$ws.Range("I7").Value = 2.06
$ws.Range("N7").Value = 3.4
$ws.Range("R7").Value = 1.3
$ws.Range("U7").Value = 1.98
$ws.Range("AF7").Value = 44
$ws.Range("AK7").Value = 400

# Row 8
$ws.Range("J8").Value = 7.8
$ws.Range("K8").Value = 8.8
$ws.Range("P8").Value = 2.66

# Row 9
$ws.Range("F9").Value = 2.74
$ws.Range("G9").Value = 2.76
$ws.Range("H9").Value = 2.7
$ws.Range("I9").Value = 2.74
$ws.Range("J9").Value = 3.7
$ws.Range("K9").Value = 3.8
$ws.Range("O9").Value = 1.27
$ws.Range("P9").Value = 2.16
$ws.Range("Q9").Value = 1.81
$ws.Range("S9").Value = 2.98
$ws.Range("T9").Value = 1.66
$ws.Range("U9").Value = 2.38
$ws.Range("X9").Value = 18
$ws.Range("Y9").Value = 13.5
$ws.Range("Z9").Value = 18.5
$ws.Range("AA9").Value = 40
$ws.Range("AB9").Value = 13
$ws.Range("AC9").Value = 8.6
$ws.Range("AD9").Value = 13
$ws.Range("AE9").Value = 28
$ws.Range("AF9").Value = 20
$ws.Range("AG9").Value = 12.5
$ws.Range("AH9").Value = 16.5
$ws.Range("AI9").Value = 38
$ws.Range("AJ9").Value = 42
$ws.Range("AK9").Value = 29
$ws.Range("AL9").Value = 75
$ws.Range("AM9").Value = 85
$ws.Range("AN9").Value = 22
$ws.Range("AO9").Value = 20

# Row 10
$ws.Range("F10").Value = 3.8
$ws.Range("G10").Value = 3.95
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 2.04
$ws.Range("J10").Value = 4
$ws.Range("P10").Value = 2.14
$ws.Range("T10").Value = 1.71
$ws.Range("U10").Value = 2.28
$ws.Range("Y10").Value = 11
$ws.Range("AA10").Value = 25
$ws.Range("AD10").Value = 10.5
$ws.Range("AI10").Value = 36
$ws.Range("AK10").Value = 46

# Row 11
$ws.Range("H11").Value = 2.28
$ws.Range("S11").Value = 2.8
$ws.Range("T11").Value = 1.64
$ws.Range("U11").Value = 2.4
$ws.Range("Z11").Value = 17
$ws.Range("AD11").Value = 12
$ws.Range("AE11").Value = 24
$ws.Range("AF11").Value = 26
$ws.Range("AG11").Value = 15
$ws.Range("AH11").Value = 16.5
$ws.Range("AI11").Value = 34
$ws.Range("AJ11").Value = 60
$ws.Range("AK11").Value = 36
$ws.Range("AL11").Value = 46
$ws.Range("AN11").Value = 28
$ws.Range("AO11").Value = 14.5

# Row 12
$ws.Range("G12").Value = 1.59
$ws.Range("N12").Value = 4.5
$ws.Range("Q12").Value = 1.77
$ws.Range("R12").Value = 1.47
$ws.Range("S12").Value = 2.88
$ws.Range("X12").Value = 21
$ws.Range("AA12").Value = 1000
$ws.Range("AK12").Value = 19
$ws.Range("AM12").Value = 1000

# Row 13
$ws.Range("F13").Value = 2.68
$ws.Range("G13").Value = 2.76
$ws.Range("I13").Value = 2.78
$ws.Range("S13").Value = 2.96
$ws.Range("T13").Value = 1.65
$ws.Range("U13").Value = 2.38

# Row 14
$ws.Range("H14").Value = 1.77
$ws.Range("K14").Value = 4.5
$ws.Range("N14").Value = 4.8
$ws.Range("X14").Value = 25
$ws.Range("Y14").Value = 13.5
$ws.Range("AA14").Value = 21
$ws.Range("AB14").Value = 25
$ws.Range("AC14").Value = 10.5
$ws.Range("AE14").Value = 18.5
$ws.Range("AG14").Value = 22
$ws.Range("AH14").Value = 20
$ws.Range("AI14").Value = 34
$ws.Range("AL14").Value = 55
$ws.Range("AM14").Value = 80
$ws.Range("AN14").Value = 48

# Row 15
$ws.Range("Q15").Value = 1.5
$ws.Range("AE15").Value = 90

# Row 16
$ws.Range("G16").Value = 5
$ws.Range("I16").Value = 2.58
$ws.Range("J16").Value = 2.78
$ws.Range("Q16").Value = 2.86

# Row 17
$ws.Range("S17").Value = 3.2
$ws.Range("AE17").Value = 110

# Row 18
$ws.Range("G18").Value = 2.02
$ws.Range("S18").Value = 3
$ws.Range("AN18").Value = 14

# Row 19
$ws.Range("N19").Value = 3.9
$ws.Range("O19").Value = 1.32
$ws.Range("P19").Value = 1.99
$ws.Range("Q19").Value = 1.95
$ws.Range("S19").Value = 3.4
$ws.Range("T19").Value = 1.75
$ws.Range("X19").Value = 15
$ws.Range("Z19").Value = 22
$ws.Range("AE19").Value = 36
$ws.Range("AF19").Value = 17
$ws.Range("AK19").Value = 27
$ws.Range("AO19").Value = 32

# Row 20
$ws.Range("F20").Value = 2.82
$ws.Range("G20").Value = 2.9
$ws.Range("H20").Value = 2.62
$ws.Range("I20").Value = 2.7
$ws.Range("Q20").Value = 1.9
$ws.Range("T20").Value = 1.73
$ws.Range("Y20").Value = 13
$ws.Range("AA20").Value = 42
$ws.Range("AB20").Value = 13.5
$ws.Range("AF20").Value = 21
$ws.Range("AG20").Value = 13.5
$ws.Range("AJ20").Value = 46
$ws.Range("AL20").Value = 42
$ws.Range("AM20").Value = 85

# Row 21
$ws.Range("F21").Value = 2.38
$ws.Range("G21").Value = 2.46
$ws.Range("H21").Value = 3.25
$ws.Range("Q21").Value = 1.9
$ws.Range("S21").Value = 3.3
$ws.Range("AC21").Value = 9.2

# Row 22
$ws.Range("Q22").Value = 1.67

# Row 23
$ws.Range("G23").Value = 2.06
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 3.85
$ws.Range("N23").Value = 4.5
$ws.Range("AM23").Value = 95

# Row 24
$ws.Range("F24").Value = 2.84
$ws.Range("AO24").Value = 27

# Row 25
$ws.Range("G25").Value = 2.58
$ws.Range("H25").Value = 2.78
$ws.Range("J25").Value = 3.9
$ws.Range("R25").Value = 1.6
$ws.Range("AH25").Value = 16.5
